# Bugfixes and extra datasets.
# Relabel the reaction-coordinate header row: the old sheet only tracked
# TS1/TS2/TS3/Product deltas; the corrected sheet reports every
# intermediate/TS free energy (relative to the reference state, RRS) for
# species 2, TS1, 3, TS2, 4 and 5, while keeping the existing TS3/Product
# columns (I/J) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 3).Value = "ΔGRRS(2)"
$ws.Cells.Item(1, 4).Value = "ΔGRRS(TS1)"
$ws.Cells.Item(1, 5).Value = "ΔGRRS(3)"
$ws.Cells.Item(1, 6).Value = "ΔGRRS(TS2)"
$ws.Cells.Item(1, 7).Value = "ΔGRRS(4)"
$ws.Cells.Item(1, 8).Value = "ΔGRRS(5)"

# Move the active selection to I1, matching the refreshed sheet view.
$null = $ws.Range("I1").Select()
